# Update the PSSM score matrix (B2:K21) with recalculated values
# from the supplemental-figures re-run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-19.3262956151879, 2.453831200291055, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, 2.466530019854751, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, 2.165499265750841, 2.83876869758963, -19.3262956151879, 2.521451856743186, -19.3262956151879, 1.880817037302376, -19.3262956151879, 2.593030746769906, -19.3262956151879),
    @(-19.3262956151879, 0.9468727928573785, -19.3262956151879, -19.3262956151879, -19.3262956151879, 1.535867244941981, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879),
    @(2.945903867038961, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, -19.3262956151879, 2.85338906864534, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879),
    @(3.619928360742919, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, 1.498340770708977, -19.3262956151879, 4.321926009898339),
    @(-19.3262956151879, -19.3262956151879, -19.3262956151879, 1.997470009907671, -19.3262956151879, 2.565147320829733, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, -19.3262956151879, 1.672320292645947, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, 2.277273232509163, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, 1.690331111141652, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, -0.2580605279926834, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, 2.210435852012807, -19.3262956151879),
    @(-19.3262956151879, 0.7085637843354706, -0.06989110191457713, -19.3262956151879, -19.3262956151879, -19.3262956151879, 0.3751817640567327, 0.8527156405906235, 1.189986328005131, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, -19.3262956151879, 0.3657766624133492, 0.9018196369726977, 1.143434627047598, -19.3262956151879),
    @(-19.3262956151879, -19.3262956151879, 1.753530931752646, -19.3262956151879, -19.3262956151879, -19.3262956151879, 1.819970370895392, 2.062584089633613, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, 1.662827893731046, 2.156341513112641, -19.3262956151879, 3.833737813950127, -19.3262956151879, 2.158604455814713, 1.924486807781014, -19.3262956151879, -19.3262956151879),
    @(-19.3262956151879, 1.72545300845031, -19.3262956151879, 2.483559241526921, -19.3262956151879, 3.483131241186148, 2.520044491428947, -19.3262956151879, -19.3262956151879, -19.3262956151879)
)

$numRows = $data.Count
$numCols = $data[0].Count
$arr = New-Object 'object[,]' $numRows,$numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r, $c] = $data[$r][$c]
    }
}

$ws.Range("B2:K21").Value = $arr
